# Auto-generated edit script: refresh LeveProfit/price columns (H-N)
# per scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 254.44444
$ws.Range("I2").Value = 290.06668
$ws.Range("J2").Value = 76.333336
$ws.Range("K2").Value = 290.06668
$ws.Range("L2").Value = 76.333336
$ws.Range("M2").Value = -177.06668
$ws.Range("N2").Value = -302.333336
$ws.Range("H28").Value = 3686.6128
$ws.Range("I28").Value = 4575.7085
$ws.Range("K28").Value = 4575.7085
$ws.Range("M28").Value = -4090.7085
$ws.Range("H33").Value = 832.6070999999999
$ws.Range("I33").Value = 145.61111
$ws.Range("J33").Value = 2069.2
$ws.Range("K33").Value = 145.61111
$ws.Range("L33").Value = 2069.2
$ws.Range("M33").Value = 83.38889
$ws.Range("N33").Value = -2527.2
$ws.Range("H42").Value = 58824024
$ws.Range("I42").Value = 100000140
$ws.Range("J42").Value = 997.8570999999999
$ws.Range("K42").Value = 300000420
$ws.Range("L42").Value = 2993.5713
$ws.Range("M42").Value = -300000190
$ws.Range("N42").Value = -3453.5713
$ws.Range("H76").Value = 3275.423
$ws.Range("I76").Value = 3089.9
$ws.Range("K76").Value = 3089.9
$ws.Range("M76").Value = -2774.9
$ws.Range("H79").Value = 3275.423
$ws.Range("I79").Value = 3089.9
$ws.Range("K79").Value = 3089.9
$ws.Range("M79").Value = -1997.9
$ws.Range("H86").Value = 4477.6665
$ws.Range("I86").Value = 4357.143
$ws.Range("J86").Value = 4899.5
$ws.Range("K86").Value = 4357.143
$ws.Range("L86").Value = 4899.5
$ws.Range("M86").Value = -3234.143
$ws.Range("N86").Value = -7145.5
$ws.Range("H89").Value = 4477.6665
$ws.Range("I89").Value = 4357.143
$ws.Range("J89").Value = 4899.5
$ws.Range("K89").Value = 21785.715
$ws.Range("L89").Value = 24497.5
$ws.Range("M89").Value = -16169.715
$ws.Range("N89").Value = -35729.5
$ws.Range("H98").Value = 9823.096
$ws.Range("J98").Value = 4128.6
$ws.Range("L98").Value = 4128.6
$ws.Range("N98").Value = -7124.6
$ws.Range("H107").Value = 338.625
$ws.Range("I107").Value = 233.13637
$ws.Range("K107").Value = 233.13637
$ws.Range("M107").Value = 1686.86363
$ws.Range("H122").Value = 9823.096
$ws.Range("J122").Value = 4128.6
$ws.Range("L122").Value = 12385.8
$ws.Range("N122").Value = -17285.8
$ws.Range("H135").Value = 1770.8125
$ws.Range("I135").Value = 1831.2174
$ws.Range("J135").Value = 1616.4445
$ws.Range("K135").Value = 16480.9566
$ws.Range("L135").Value = 14548.0005
$ws.Range("M135").Value = -13945.9566
$ws.Range("N135").Value = -19618.0005
$ws.Range("H137").Value = 19982.166
$ws.Range("I137").Value = 22978.6
$ws.Range("K137").Value = 68935.79999999999
$ws.Range("M137").Value = -66385.79999999999
$ws.Range("H141").Value = 5009.5
$ws.Range("I141").Value = 2568.6843
$ws.Range("J141").Value = 8576.846
$ws.Range("K141").Value = 7706.0529
$ws.Range("L141").Value = 25730.538
$ws.Range("M141").Value = -2526.0529
$ws.Range("N141").Value = -36090.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4589.4443
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H14").Value = 1014.1429
$ws.Range("I14").Value = 634
$ws.Range("J14").Value = 1299.25
$ws.Range("K14").Value = 634
$ws.Range("L14").Value = 1299.25
$ws.Range("M14").Value = -459
$ws.Range("N14").Value = -1649.25
$ws.Range("H32").Value = 4718
$ws.Range("I32").Value = 4938
$ws.Range("K32").Value = 4938
$ws.Range("M32").Value = -4651
$ws.Range("H38").Value = 4456.1665
$ws.Range("I38").Value = 4163.2
$ws.Range("K38").Value = 4163.2
$ws.Range("M38").Value = -3696.2
$ws.Range("H61").Value = 4943.8696
$ws.Range("I61").Value = 4532.227
$ws.Range("J61").Value = 14000
$ws.Range("K61").Value = 4532.227
$ws.Range("L61").Value = 14000
$ws.Range("M61").Value = -4320.227
$ws.Range("N61").Value = -14424
$ws.Range("H116").Value = 4589.4443
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("H122").Value = 3018.2307
$ws.Range("I122").Value = 3075.4
$ws.Range("J122").Value = 2982.5
$ws.Range("K122").Value = 9226.200000000001
$ws.Range("L122").Value = 8947.5
$ws.Range("M122").Value = -6776.200000000001
$ws.Range("N122").Value = -13847.5
$ws.Range("H132").Value = 2694.875
$ws.Range("I132").Value = 1968.1666
$ws.Range("J132").Value = 4875
$ws.Range("K132").Value = 5904.4998
$ws.Range("L132").Value = 14625
$ws.Range("M132").Value = -3374.4998
$ws.Range("N132").Value = -19685
$ws.Range("H136").Value = 4943.8696
$ws.Range("I136").Value = 4532.227
$ws.Range("J136").Value = 14000
$ws.Range("K136").Value = 13596.681
$ws.Range("L136").Value = 42000
$ws.Range("M136").Value = -11046.681
$ws.Range("N136").Value = -47100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4589.4443
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("H7").Value = 975
$ws.Range("J7").Value = 1450
$ws.Range("L7").Value = 1450
$ws.Range("N7").Value = -1676
$ws.Range("H12").Value = 94.75
$ws.Range("I12").Value = 43.333332
$ws.Range("J12").Value = 249
$ws.Range("K12").Value = 43.333332
$ws.Range("L12").Value = 249
$ws.Range("M12").Value = 124.666668
$ws.Range("N12").Value = -585
$ws.Range("H22").Value = 1179.2
$ws.Range("I22").Value = 958.9
$ws.Range("K22").Value = 958.9
$ws.Range("M22").Value = -785.9
$ws.Range("H61").Value = 74998
$ws.Range("J61").Value = 74998
$ws.Range("L61").Value = 74998
$ws.Range("N61").Value = -75624
$ws.Range("H105").Value = 2893.5715
$ws.Range("I105").Value = 2095
$ws.Range("K105").Value = 2095
$ws.Range("M105").Value = -348
$ws.Range("H134").Value = 2418.0908
$ws.Range("I134").Value = 2326.2104
$ws.Range("K134").Value = 6978.6312
$ws.Range("M134").Value = -4443.6312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 848.75
$ws.Range("J19").Value = 2710
$ws.Range("L19").Value = 2710
$ws.Range("N19").Value = -3050
$ws.Range("H24").Value = 848.75
$ws.Range("J24").Value = 2710
$ws.Range("L24").Value = 2710
$ws.Range("N24").Value = -3050
$ws.Range("H31").Value = 5041
$ws.Range("I31").Value = 2101.6428
$ws.Range("K31").Value = 2101.6428
$ws.Range("M31").Value = -1806.6428
$ws.Range("H34").Value = 5041
$ws.Range("I34").Value = 2101.6428
$ws.Range("K34").Value = 2101.6428
$ws.Range("M34").Value = -1899.6428
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30368
$ws.Range("H58").Value = 3212.5
$ws.Range("I58").Value = 3212.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3212.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = None
$ws.Range("N58").Value = ""
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H132").Value = 5031.25
$ws.Range("I132").Value = 5041.6665
$ws.Range("K132").Value = 15124.9995
$ws.Range("M132").Value = -12594.9995
$ws.Range("H135").Value = 249999
$ws.Range("J135").Value = 249999
$ws.Range("L135").Value = 249999
$ws.Range("N135").Value = -260139
$ws.Range("H136").Value = 3212.5
$ws.Range("I136").Value = 3212.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9637.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = None
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2505.7693
$ws.Range("J23").Value = 2020.1666
$ws.Range("L23").Value = 6060.4998
$ws.Range("N23").Value = -6530.4998
$ws.Range("H34").Value = 1485.3704
$ws.Range("J34").Value = 1485.3704
$ws.Range("L34").Value = 4456.1112
$ws.Range("N34").Value = -4624.1112
$ws.Range("H55").Value = 364411.84
$ws.Range("I55").Value = 833776.8
$ws.Range("J55").Value = 12388.125
$ws.Range("K55").Value = 2501330.4
$ws.Range("L55").Value = 37164.375
$ws.Range("M55").Value = -2501153.4
$ws.Range("N55").Value = -37518.375
$ws.Range("H63").Value = 999
$ws.Range("I63").Value = 999
$ws.Range("K63").Value = 2997
$ws.Range("M63").Value = -2248
$ws.Range("H66").Value = 999
$ws.Range("I66").Value = 999
$ws.Range("K66").Value = 8991
$ws.Range("M66").Value = -5247
$ws.Range("H88").Value = 14181.818
$ws.Range("J88").Value = 19428.572
$ws.Range("L88").Value = 58285.716
$ws.Range("N88").Value = -59141.716
$ws.Range("H91").Value = 14181.818
$ws.Range("J91").Value = 19428.572
$ws.Range("L91").Value = 58285.716
$ws.Range("N91").Value = -61249.716
$ws.Range("H132").Value = 2266.9412
$ws.Range("I132").Value = 1266.1428
$ws.Range("J132").Value = 2967.5
$ws.Range("K132").Value = 11395.2852
$ws.Range("L132").Value = 26707.5
$ws.Range("M132").Value = -8865.2852
$ws.Range("N132").Value = -31767.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 518
$ws.Range("I9").Value = 518
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 518
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = None
$ws.Range("N9").Value = ""
$ws.Range("H49").Value = 25798.2
$ws.Range("J49").Value = 25798.2
$ws.Range("L49").Value = 25798.2
$ws.Range("N49").Value = -26166.2
$ws.Range("H52").Value = 272500
$ws.Range("J52").Value = 272500
$ws.Range("L52").Value = 272500
$ws.Range("N52").Value = -273018
$ws.Range("H70").Value = 10701.209
$ws.Range("I70").Value = 11751.5
$ws.Range("J70").Value = 10634.523
$ws.Range("K70").Value = 11751.5
$ws.Range("L70").Value = 10634.523
$ws.Range("M70").Value = -11481.5
$ws.Range("N70").Value = -11174.523
$ws.Range("H73").Value = 10701.209
$ws.Range("I73").Value = 11751.5
$ws.Range("J73").Value = 10634.523
$ws.Range("K73").Value = 11751.5
$ws.Range("L73").Value = 10634.523
$ws.Range("M73").Value = -10815.5
$ws.Range("N73").Value = -12506.523
$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884
$ws.Range("H122").Value = 9455.654
$ws.Range("I122").Value = 9455.654
$ws.Range("K122").Value = 28366.962
$ws.Range("M122").Value = -25916.962

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2013.5
$ws.Range("I22").Value = 1989.5454
$ws.Range("J22").Value = 2026.675
$ws.Range("K22").Value = 1989.5454
$ws.Range("L22").Value = 2026.675
$ws.Range("M22").Value = -1694.5454
$ws.Range("N22").Value = -2616.675
$ws.Range("H27").Value = 2013.5
$ws.Range("I27").Value = 1989.5454
$ws.Range("J27").Value = 2026.675
$ws.Range("K27").Value = 1989.5454
$ws.Range("L27").Value = 2026.675
$ws.Range("M27").Value = -1882.5454
$ws.Range("N27").Value = -2240.675
$ws.Range("H55").Value = 1832.2858
$ws.Range("I55").Value = 201.57143
$ws.Range("J55").Value = 3463
$ws.Range("K55").Value = 201.57143
$ws.Range("L55").Value = 3463
$ws.Range("M55").Value = -28.57142999999999
$ws.Range("N55").Value = -3809
$ws.Range("H61").Value = 3422.9092
$ws.Range("I61").Value = 1105.3077
$ws.Range("K61").Value = 1105.3077
$ws.Range("M61").Value = -903.3077000000001
$ws.Range("H113").Value = 3422.9092
$ws.Range("I113").Value = 1105.3077
$ws.Range("K113").Value = 1105.3077
$ws.Range("M113").Value = 1064.6923
$ws.Range("H122").Value = 3427.2856
$ws.Range("I122").Value = 3581.8333
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 10745.4999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -8295.499899999999
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 3605.027
$ws.Range("I132").Value = 3282.9
$ws.Range("J132").Value = 4985.5713
$ws.Range("K132").Value = 9848.700000000001
$ws.Range("L132").Value = 14956.7139
$ws.Range("M132").Value = -7318.700000000001
$ws.Range("N132").Value = -20016.7139
$ws.Range("H134").Value = 68000
$ws.Range("J134").Value = 68000
$ws.Range("L134").Value = 68000
$ws.Range("N134").Value = -78140
$ws.Range("H136").Value = 2998.4138
$ws.Range("I136").Value = 3102.2
$ws.Range("J136").Value = 2349.75
$ws.Range("K136").Value = 9306.599999999999
$ws.Range("L136").Value = 7049.25
$ws.Range("M136").Value = -6756.599999999999
$ws.Range("N136").Value = -12149.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 6666.6665
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = -4860
$ws.Range("N8").Value = -10280
$ws.Range("H132").Value = 8139.25
$ws.Range("I132").Value = 8139.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 24417.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = None
$ws.Range("N132").Value = ""
